$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates from the commit diff.
# Numeric-looking text values are forced to stay text (matching the
# original inlineStr/shared-string cell type) by temporarily applying a
# "@" (Text) number format before the assignment, then resetting the
# cell style back to "Normal" so no lingering style index is left behind.

$ws.Range("D2").Value = "24.908.99"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "1.704.74"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4012"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4070"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.470"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08829"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.497"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.047"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001351"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "1.634.53"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07192"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "21.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.246"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("D24").Value = "24.891.91"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.330"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.890"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.485"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +23.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "144.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.212"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.268"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.60%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08755"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.396"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("B35").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C35").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D35").Value = "1.815.71"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03188"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.037"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8521"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09446"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.476"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.721"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7471"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.230"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.397"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08401"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.13%  "
